$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts everything below down by one).
$ws.Rows("3:3").Insert()

# New row 3 content: Roberto / 0871 / Casa Alexandre / Duas câmeras sem imagem.
$ws.Range("A3").Value = "Roberto"
$ws.Range("B3").Value = "'0871"
$ws.Range("C3").Value = "Casa Alexandre"
$ws.Range("D3").Value = "Duas câmeras sem imagem."

# Copy formatting from row 2 onto the new row 3 so it matches (border/fill/alignment/numFmt).
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in follow-up/observation info on the (now shifted) rows 6 and 7.
$ws.Range("F6").Value = "Fonte de 10AP foi trocada."
$ws.Range("G6").Value = "Concluido"
$ws.Range("E6").Value = "Comunicação foi restaurada e passamos pro DDNS."

$ws.Range("E7").Value = "Local em reforma, foi combinado com Edenis dele retornar quando finalizar."
$ws.Range("G7").Value = "Concluido"

# Reassign technician for rows 8-10 (now shifted) from Giovani to Marcos.
$ws.Range("A8").Value = "Marcos"
$ws.Range("A9").Value = "Marcos"
$ws.Range("A10").Value = "Marcos"

# Populate the previously-empty row 15 (now shifted) with a new agenda entry.
$ws.Range("A15").Value = "Fábio"
$ws.Range("B15").Value = "'0730"
$ws.Range("C15").Value = "JR Leo Ipatinga"
$ws.Range("D15").Value = "Central de alarmes sem comunicação."

# Update the sheet view: scroll back to column A and move the active selection to D11.
$ws.Range("D11").Select()
